$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before S (pushes old S:AD to U:AF)
$ws.Range("S1:T1").EntireColumn.Insert()

# Set new header labels (set T1 first so shared-string table order matches: ci_ditch then co_ditch)
$ws.Range("T1").Value2 = "ci_ditch"
$ws.Range("S1").Value2 = "co_ditch"

# Update changed values in row 2 (using new column letters, post-insert)
$ws.Range("C2").Value2 = 1
$ws.Range("E2").Value2 = -0.001
$ws.Range("H2").Value2 = 0.2
$ws.Range("J2").Value2 = 0.001
$ws.Range("K2").Value2 = 5
$ws.Range("L2").Value2 = 2.5
$ws.Range("P2").Value2 = 2.2000000000000002
$ws.Range("S2").Value2 = 0.5
$ws.Range("T2").Value2 = 1
$ws.Range("U2").Value2 = 0.75
$ws.Range("V2").Value2 = 0
$ws.Range("Y2").Value2 = 2

# Update the active cell selection to match the saved view state
$ws.Range("O2").Select()

Write-Host "done"
